$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '2021-06-30'
$ws.Range("H2").Value = 'Madrigal'
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 21000
$ws.Range("M2").Value = 20333
$ws.Range("N2").Value = '$/caja 40 unidades'
$ws.Range("O2").Value = 'Región de Coquimbo'
$ws.Range("P2").Value = 508
$ws.Range("Q2").Value = 40

# Row 3
$ws.Range("D3").Value = '2021-06-30'
$ws.Range("H3").Value = 'Symphony'
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 21000
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 21500
$ws.Range("N3").Value = '$/caja 40 unidades'
$ws.Range("P3").Value = 538
$ws.Range("Q3").Value = 40

# Row 4
$ws.Range("D4").Value = '2021-06-09'
$ws.Range("H4").Value = 'Argentina(o)'
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19500
$ws.Range("N4").Value = '$/caja 50 unidades'
$ws.Range("P4").Value = 390
$ws.Range("Q4").Value = 50

# Row 5
$ws.Range("D5").Value = '2020-11-25'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("P5").Value = 362
$ws.Range("Q5").Value = 40

# Row 6
$ws.Range("D6").Value = '2021-07-28'
$ws.Range("H6").Value = 'Madrigal'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 21000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21500
$ws.Range("P6").Value = 538

# Row 8
$ws.Range("D8").Value = '2021-08-11'
$ws.Range("H8").Value = 'Symphony'
$ws.Range("K8").Value = 21000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21500
$ws.Range("N8").Value = '$/caja 50 unidades'
$ws.Range("P8").Value = 430
$ws.Range("Q8").Value = 50

# Row 9
$ws.Range("D9").Value = '2021-06-23'
$ws.Range("H9").Value = 'Argentina(o)'
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 21000
$ws.Range("M9").Value = 20429
$ws.Range("N9").Value = '$/caja 50 unidades'
$ws.Range("P9").Value = 409
$ws.Range("Q9").Value = 50

# Row 10
$ws.Range("D10").Value = '2021-06-23'
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 22000
$ws.Range("L10").Value = 23000
$ws.Range("M10").Value = 22500
$ws.Range("P10").Value = 562

# Row 11
$ws.Range("D11").Value = '2021-07-07'
$ws.Range("H11").Value = 'Madrigal'
$ws.Range("J11").Value = 80

# Row 12
$ws.Range("D12").Value = '2021-07-07'
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 30
$ws.Range("M12").Value = 19333
$ws.Range("N12").Value = '$/caja 50 unidades'
$ws.Range("P12").Value = 387
$ws.Range("Q12").Value = 50

# Row 13
$ws.Range("D13").Value = '2021-07-07'
$ws.Range("H13").Value = 'Symphony'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 21000
$ws.Range("M13").Value = 20400
$ws.Range("P13").Value = 510

# Row 14
$ws.Range("D14").Value = '2021-08-25'
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 19000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19500
$ws.Range("P14").Value = 488

# Row 15
$ws.Range("D15").Value = '2021-08-18'
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("P15").Value = 488
$ws.Range("Q15").Value = 40

# Row 16
$ws.Range("D16").Value = '2021-07-21'
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = 21000
$ws.Range("L16").Value = 22000
$ws.Range("M16").Value = 21500
$ws.Range("P16").Value = 538

# Row 17
$ws.Range("D17").Value = '2021-07-14'
$ws.Range("J17").Value = 140
$ws.Range("K17").Value = 21000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21500
$ws.Range("P17").Value = 538

# Row 18
$ws.Range("D18").Value = '2021-08-27'
$ws.Range("H18").Value = 'Madrigal'
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19500
$ws.Range("N18").Value = '$/caja 40 unidades'
$ws.Range("P18").Value = 488
$ws.Range("Q18").Value = 40

# Row 19
$ws.Range("D19").Value = '2021-08-04'
$ws.Range("H19").Value = 'Symphony'
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = 21000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 21500
$ws.Range("P19").Value = 538

# Row 20
$ws.Range("D20").Value = '2020-12-02'
$ws.Range("H20").Value = 'Española'
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13500
$ws.Range("N20").Value = '$/caja 30 unidades'
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 450
$ws.Range("Q20").Value = 30
